# Rename the "library construction" sheet and touch up the related
# bits of the workbook (table name label, column width, selections).

$wb = $excel.ActiveWorkbook

$isaSheet = $wb.Worksheets.Item(1)
$libSheet = $wb.Worksheets.Item(2)

# 1. Rename the second sheet from "single_cell_library_constructio"
#    (Excel's 31-char sheet-name truncation of the original) to the
#    shorter "single_cell_library".
$libSheet.Name = "single_cell_library"

# 2. The "Table" metadata row on the isa_template sheet (B7) stores the
#    same name as plain text - keep it in sync.
$isaSheet.Range("B7").Value = "single_cell_library"

# 3. Widen column A on the isa_template sheet so the longer labels fit.
$isaSheet.Columns.Item(1).ColumnWidth = 41.616666666666667

# 4. Restore the selection/active-cell bookkeeping that Excel records
#    per sheet: isa_template now has B8 selected, and the library sheet
#    is back to the default A1 selection (no longer carrying the old
#    F13:F17 remnant).
$isaSheet.Activate() | Out-Null
$isaSheet.Range("B8").Select() | Out-Null

$libSheet.Activate() | Out-Null
$libSheet.Range("A1").Select() | Out-Null

$isaSheet.Activate() | Out-Null
